# add PI Loop Filter
# Update the "Branch LPF" sheet cutoff frequency input (B2): 175 -> 600
# Update the "LFz" sheet cutoff frequency input (B2): 480 -> 100
# Update the "LFz" sheet scale-bits input (B18): 15 -> 13

$wb = $excel.ActiveWorkbook

$branchLpf = $wb.Worksheets.Item("Branch LPF")
$branchLpf.Range("B2").Value = 600
$branchLpf.Activate()
[void]$branchLpf.Range("B18").Select()

$lfz = $wb.Worksheets.Item("LFz")
$lfz.Range("B2").Value = 100
$lfz.Range("B18").Value = 13

# Mark the LFz sheet as the active / selected sheet, matching the
# final tabSelected state in the workbook.
$lfz.Activate()
[void]$lfz.Range("B19").Select()

$wb.Save()
